$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for the week of 25_02_2024
$ws.Range("D1").Value = "25_02_2024"

# Add the new appointment counts per recepcionista
$ws.Range("D2").Value = 1245
$ws.Range("D3").Value = 1281
$ws.Range("D4").Value = 1275
$ws.Range("D5").Value = 3015
$ws.Range("D6").Value = 83

# Update the selection to match the post-edit state
$ws.Range("D3").Select()
